# Update evidence for dee-effes

$wb = $excel.ActiveWorkbook

$wsInfo = $wb.Worksheets.Item("Info")
$wsA1   = $wb.Worksheets.Item("A1")
$wsA2   = $wb.Worksheets.Item("A2")
$wsA3   = $wb.Worksheets.Item("A3")
$wsA4   = $wb.Worksheets.Item("A4")
$wsA5   = $wb.Worksheets.Item("A5")
$wsA6   = $wb.Worksheets.Item("A6")

# ---- A1 sheet: update TxHash / ClassID values ----
$wsA1.Range("A2").Value = "444971A9CE8D7F91F05F790791C437569F7920C73B0860D7FA481A53DC26B6CB"
$wsA1.Range("B2").Value = "blcsBenni001"

# ---- A2 sheet: update TxHash / ClassID / NFTID values ----
$wsA2.Range("A2").Value = "81C700D246F57039F3BDC16B0DA0C7CFF937356271269FD8CF50A87316923F9D"
$wsA2.Range("B2").Value = "blcsBenni001"
$wsA2.Range("C2").Value = "blcsNFTBenni002"
$wsA2.Range("A3").Value = "4AED26214C9A72F97E381AEB3385199F14BDA6324448FAB0029202D092C80018"
$wsA2.Range("B3").Value = "blcsBenni001"
$wsA2.Range("C3").Value = "blcsNFTmars001"
$wsA2.Columns.Item(2).ColumnWidth = 14

# ---- A3 sheet: update TxHash / ClassID values ----
$wsA3.Range("A2").Value = "5F0E44A16797DD129730B38EFA84F240C3D8EA8A82384DBBB7BFBDCA2AA4B916"
$wsA3.Range("B2").Value = "stars1yxv9njz6nfpwxtxzfawe5a2mv7z9gu9hl23zhu7g5c4jaxh7y4nqjkav4n"

# ---- A4 sheet: add row 2 ----
$wsA4.Range("A2").Value = "1A0C443A20A619C1DE7858139C0B13FA77A4CFA463D7AF035C67A65B5FD318E1"
$wsA4.Range("B2").Value = "ibc/6EAF424647FC10605DA6FFF50CF3FD26D3AA49495FA211FE0551F84DD9E79C7F"
$wsA4.Range("C2").Value = "blcsNFTBenni005"
$wsA4.Range("D2").Value = "uptick_7000-2"
$wsA4.Columns.Item(1).ColumnWidth = 71.42578125
$wsA4.Columns.Item(2).ColumnWidth = 73.7109375
$wsA4.Columns.Item(3).ColumnWidth = 16.140625
$wsA4.Columns.Item(4).ColumnWidth = 13.28515625

# ---- A5 sheet: add row 2 ----
$wsA5.Range("A2").Value = "B1FC2856C78D496EE0F4541EA15FCE2DAC5D630E2C8047C59A27CB02726241BA"
$wsA5.Range("B2").Value = "stars1yxv9njz6nfpwxtxzfawe5a2mv7z9gu9hl23zhu7g5c4jaxh7y4nqjkav4n"
$wsA5.Range("C2").Value = "blcsNFTBenni002"
$wsA5.Range("D2").Value = "elgafar-1"
$wsA5.Columns.Item(1).ColumnWidth = 72
$wsA5.Columns.Item(2).ColumnWidth = 67.5703125
$wsA5.Columns.Item(3).ColumnWidth = 17

# ---- A6 sheet: add row 2 ----
$wsA6.Range("A2").Value = "EA84CC14BE7E279EC05193B3FDD6FBBBB55702D16BF654C7EBA6E6D6CF7AA736"
$wsA6.Range("B2").Value = "ibc/6EAF424647FC10605DA6FFF50CF3FD26D3AA49495FA211FE0551F84DD9E79C7F"
$wsA6.Range("C2").Value = "blcsNFTBenni005"
$wsA6.Range("D2").Value = "uptick_7000-2"
$wsA6.Columns.Item(1).ColumnWidth = 73
$wsA6.Columns.Item(2).ColumnWidth = 73.140625
$wsA6.Columns.Item(3).ColumnWidth = 15.7109375
$wsA6.Columns.Item(4).ColumnWidth = 13

# ---- Selections / view state ----
$wsInfo.Range("G6").Select()
$wsA1.Range("A5").Select()
$wsA2.Range("A1:D1").Select()
$wsA3.Range("B9").Select()
$wsA4.Range("B2").Select()
$wsA5.Range("A1:D1").Select()
$wsA5.Range("D1").Activate()
$wsA6.Range("D2").Select()

# Make A6 the active / selected sheet and leave Info no longer tab-selected
$wsA6.Activate()
